$d = $word.ActiveDocument

# Locate the very last paragraph in the document body (the one holding the
# inline image) and split right after it so the image paragraph itself is
# left completely untouched.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

# The freshly created (empty) paragraph is now the new last paragraph;
# collapse its range to its start and replace/fill it via InsertXML with
# the six paragraphs described by the diff.
$afterIndex = $d.Paragraphs.Count
$insertionPara = $d.Paragraphs.Item($afterIndex)
$insertionRange = $insertionPara.Range
$insertionRange.Collapse(1)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:jc w:val="both"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Novo paágrafo inserido no ramo 1</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:jc w:val="both"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0" w:line="360" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="DauphinPlain" w:hAnsi="DauphinPlain"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="DauphinPlain" w:hAnsi="DauphinPlain"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>Why do we use it?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>It is a long established fact that a reader will be distracted by the readable content of a page when looking at its layout. The point of using Lorem Ipsum is that it has a more-or-less normal distribution of letters, as opposed to using 'Content here, content here', making it look like readable English. Many desktop publishing packages and web page editors now use Lorem Ipsum as their default model text, and a search for 'lorem ipsum' will uncover many web sites still in their infancy. Various versions have evolved over the years, sometimes by accident, sometimes on purpose (injected humour and the like).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:jc w:val="both"/></w:pPr></w:p>
'@

$insertionRange.InsertXML($xml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
